# Improve the localT2 values: refresh accuracy figures in column A.
# A1 keeps the header label; A2:A49 get updated numeric accuracy values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_R_acc_LT"

$ws.Range("A2").Value = 80.915331807780319
$ws.Range("A3").Value = 83.981693363844386
$ws.Range("A4").Value = 84.302059496567509
$ws.Range("A5").Value = 88.054919908466829
$ws.Range("A6").Value = 88.054919908466829
$ws.Range("A7").Value = 87.643020594965677
$ws.Range("A8").Value = 78.764302059496572
$ws.Range("A9").Value = 80.549199084668189
$ws.Range("A10").Value = 79.54233409610984
$ws.Range("A11").Value = 78.535469107551492
$ws.Range("A12").Value = 74.233409610983983
$ws.Range("A13").Value = 79.77116704805492
$ws.Range("A14").Value = 77.848970251716239
$ws.Range("A15").Value = 78.169336384439362
$ws.Range("A16").Value = 80.45766590389016
$ws.Range("A17").Value = 77.66590389016018
$ws.Range("A18").Value = 79.816933638443928
$ws.Range("A19").Value = 85.080091533180777
$ws.Range("A20").Value = 88.009153318077807
$ws.Range("A21").Value = 88.009153318077807
$ws.Range("A22").Value = 87.826086956521749
$ws.Range("A23").Value = 79.496567505720833
$ws.Range("A24").Value = 82.10526315789474
$ws.Range("A25").Value = 81.51029748283753
$ws.Range("A26").Value = 81.64759725400458
$ws.Range("A27").Value = 80.82379862700229
$ws.Range("A28").Value = 81.693363844393602
$ws.Range("A29").Value = 81.28146453089245
$ws.Range("A30").Value = 80.274599542334087
$ws.Range("A31").Value = 80.183066361556072
$ws.Range("A32").Value = 89.473684210526315
$ws.Range("A33").Value = 92.494279176201374
$ws.Range("A34").Value = 92.082379862700222
$ws.Range("A35").Value = 83.203661327231131
$ws.Range("A36").Value = 88.054919908466829
$ws.Range("A37").Value = 68.054919908466815
$ws.Range("A38").Value = 84.668192219679639
$ws.Range("A39").Value = 80.869565217391298
$ws.Range("A40").Value = 79.588100686498848
$ws.Range("A41").Value = 79.450800915331811
$ws.Range("A42").Value = 79.54233409610984
$ws.Range("A43").Value = 79.725400457665913
$ws.Range("A44").Value = 79.679633867276891
$ws.Range("A45").Value = 83.020594965675059
$ws.Range("A46").Value = 84.713958810068647
$ws.Range("A47").Value = 78.993135011441652
$ws.Range("A48").Value = 77.848970251716239
$ws.Range("A49").Value = 80.45766590389016
